# Update Name of Algo - adjust specific numeric result values on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.44099999999999

$ws.Range("A4").Value = -21.22069999999999
$ws.Range("B4").Value = 4.802100000000003
$ws.Range("C4").Value = -11.16769999999999

$ws.Range("B5").Value = 5.365899999999999

$ws.Range("A6").Value = -21.67090000000002
$ws.Range("B6").Value = 5.426699999999999

$ws.Range("A7").Value = -21.35320000000002

$ws.Range("A8").Value = -21.53050000000002
$ws.Range("B8").Value = 4.811400000000001

$ws.Range("C9").Value = -11.6619

$ws.Range("C11").Value = -13.87070000000001

$ws.Range("C14").Value = -11.65819999999999

$ws.Range("A16").Value = -21.51680000000003
$ws.Range("B16").Value = 5.439299999999995

$ws.Range("C18").Value = -14.54930000000001

$ws.Range("A20").Value = -22.85110000000002

$ws.Range("A21").Value = -20.6301

$ws.Range("B22").Value = 5.496999999999998

$ws.Range("C25").Value = -11.3876
